$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "M2DocEvaluator.java:1267" "M2DocEvaluator.java:1313"
Replace-Text "M2DocEvaluator.java:1038" "M2DocEvaluator.java:1084"
Replace-Text "M2DocEvaluator.java:1254" "M2DocEvaluator.java:1300"
Replace-Text "M2DocEvaluator.java:1278" "M2DocEvaluator.java:1324"
Replace-Text "M2DocEvaluator.java:275)" "M2DocEvaluator.java:278)"
Replace-Text "M2DocEvaluator.java:264)" "M2DocEvaluator.java:267)"
Replace-Text "M2DocUtils.java:712" "M2DocUtils.java:694"
Replace-Text "AbstractTemplatesTestSuite.java:459" "AbstractTemplatesTestSuite.java:475"
Replace-Text "AbstractTemplatesTestSuite.java:369" "AbstractTemplatesTestSuite.java:384"
